$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -10.37282494035166
$ws.Range("J25").Value = -10.51070347022909
$ws.Range("K25").Value = 3.635113502123218
$ws.Range("I26").Value = -11.11602434548657
$ws.Range("J26").Value = 3.029792626865739
$ws.Range("K26").Value = 5.056794673020033
$ws.Range("H27").Value = -10.60022236106626
$ws.Range("I27").Value = 3.545594611286049
$ws.Range("J27").Value = 5.572596657440343
$ws.Range("K27").Value = -3.017090008061245
$ws.Range("G28").Value = -9.61602434548657
$ws.Range("H28").Value = 4.529792626865739
$ws.Range("I28").Value = 6.556794673020033
$ws.Range("J28").Value = -2.032891992481555
$ws.Range("K28").Value = -5.007106464233121
$ws.Range("F29").Value = -10.72192254212553
$ws.Range("G29").Value = 3.423894430226781
$ws.Range("H29").Value = 5.450896476381075
$ws.Range("I29").Value = -3.138790189120513
$ws.Range("J29").Value = -6.113004660872079
$ws.Range("K29").Value = 1.844819121192472
$ws.Range("E30").Value = -10.38364849008269
$ws.Range("F30").Value = 3.762168482269615
$ws.Range("G30").Value = 5.789170528423909
$ws.Range("H30").Value = -2.800516137077679
$ws.Range("I30").Value = -5.774730608829245
$ws.Range("J30").Value = 2.183093173235306
$ws.Range("K30").Value = 0.567070620731144
$ws.Range("D31").Value = -10.79732299193613
$ws.Range("E31").Value = 3.34849398041618
$ws.Range("F31").Value = 5.375496026570474
$ws.Range("G31").Value = -3.214190638931115
$ws.Range("H31").Value = -6.18840511068268
$ws.Range("I31").Value = 1.76941867138187
$ws.Range("J31").Value = 0.1533961188777085
$ws.Range("K31").Value = -0.7993314525965047
$ws.Range("C32").Value = -14.01602434548657
$ws.Range("D32").Value = 0.129792626865739
$ws.Range("E32").Value = 2.156794673020033
$ws.Range("F32").Value = -6.432891992481555
$ws.Range("G32").Value = -9.40710646423312
$ws.Range("H32").Value = -1.44928268216857
$ws.Range("I32").Value = -3.065305234672732
$ws.Range("J32").Value = -4.018032806146945
$ws.Range("K32").Value = -4.237928492263277
$ws.Range("B33").Value = -21.73309111108851
$ws.Range("C33").Value = -7.5872741387362
$ws.Range("D33").Value = -5.560272092581906
$ws.Range("E33").Value = -14.14995875808349
$ws.Range("F33").Value = -17.12417322983506
$ws.Range("G33").Value = -9.166349447770509
$ws.Range("H33").Value = -10.78237200027467
$ws.Range("I33").Value = -11.73509957174888
$ws.Range("J33").Value = -11.95499525786522
$ws.Range("K33").Value = -9.913133853543295
$ws.Range("B34").Value = 14.14581697235231
$ws.Range("C34").Value = 16.1728190185066
$ws.Range("D34").Value = 7.583132353005015
$ws.Range("E34").Value = 4.608917881253449
$ws.Range("F34").Value = 12.566741663318
$ws.Range("G34").Value = 10.95071911081384
$ws.Range("H34").Value = 9.997991539339624
$ws.Range("I34").Value = 9.778095853223293
$ws.Range("J34").Value = 11.81995725754521
$ws.Range("K34").Value = 12.04588566150899
$ws.Range("B35").Value = 2.027002046154294
$ws.Range("C35").Value = -6.562684619347294
$ws.Range("D35").Value = -9.53689909109886
$ws.Range("E35").Value = -1.579075309034309
$ws.Range("F35").Value = -3.195097861538471
$ws.Range("G35").Value = -4.147825433012684
$ws.Range("H35").Value = -4.367721119129016
$ws.Range("I35").Value = -2.325859714807095
$ws.Range("J35").Value = -2.099931310843317
$ws.Range("K35").Value = -3.377147207894367
$ws.Range("B36").Value = -8.589686665501588
$ws.Range("C36").Value = -11.56390113725315
$ws.Range("D36").Value = -3.606077355188603
$ws.Range("E36").Value = -5.222099907692765
$ws.Range("F36").Value = -6.174827479166979
$ws.Range("G36").Value = -6.39472316528331
$ws.Range("H36").Value = -4.352861760961389
$ws.Range("I36").Value = -4.126933356997611
$ws.Range("J36").Value = -5.404149254048662
$ws.Range("K36").Value = -5.09111828891804
$ws.Range("B37").Value = -2.974214471751566
$ws.Range("C37").Value = 4.983609310312985
$ws.Range("D37").Value = 3.367586757808823
$ws.Range("E37").Value = 2.41485918633461
$ws.Range("F37").Value = 2.194963500218279
$ws.Range("G37").Value = 4.2368249045402
$ws.Range("H37").Value = 4.462753308503977
$ws.Range("I37").Value = 3.185537411452927
$ws.Range("J37").Value = 3.498568376583549
$ws.Range("K37").Value = 2.99615326272577
$ws.Range("B38").Value = 7.957823782064551
$ws.Range("C38").Value = 6.341801229560389
$ws.Range("D38").Value = 5.389073658086176
$ws.Range("E38").Value = 5.169177971969845
$ws.Range("F38").Value = 7.211039376291765
$ws.Range("G38").Value = 7.436967780255543
$ws.Range("H38").Value = 6.159751883204493
$ws.Range("I38").Value = 6.472782848335115
$ws.Range("J38").Value = 5.970367734477335
$ws.Range("K38").Value = 6.62200098188616
$ws.Range("B39").Value = -1.616022552504162
$ws.Range("C39").Value = -2.568750123978375
$ws.Range("D39").Value = -2.788645810094706
$ws.Range("E39").Value = -0.7467844057727859
$ws.Range("F39").Value = -0.5208560018090078
$ws.Range("G39").Value = -1.798071898860058
$ws.Range("H39").Value = -1.485040933729436
$ws.Range("I39").Value = -1.987456047587215
$ws.Range("J39").Value = -1.335822800178391
$ws.Range("K39").Value = -1.528722419220372
$ws.Range("B40").Value = -0.9527275714742132
$ws.Range("C40").Value = -1.172623257590544
$ws.Range("D40").Value = 0.8692381467313761
$ws.Range("E40").Value = 1.095166550695154
$ws.Range("F40").Value = -0.1820493463558961
$ws.Range("G40").Value = 0.1309816187747259
$ws.Range("H40").Value = -0.3714334950830534
$ws.Range("I40").Value = 0.2801997523257711
$ws.Range("J40").Value = 0.0873001332837901
$ws.Range("K40").Value = -0.4234413153011048
$ws.Range("B41").Value = -0.2198956861163308
$ws.Range("C41").Value = 1.821965718205589
$ws.Range("D41").Value = 2.047894122169367
$ws.Range("E41").Value = 0.7706782251183171
$ws.Range("F41").Value = 1.083709190248939
$ws.Range("G41").Value = 0.5812940763911598
$ws.Range("H41").Value = 1.232927323799984
$ws.Range("I41").Value = 1.040027704758003
$ws.Range("J41").Value = 0.5292862561731084
$ws.Range("K41").Value = 0.7346869887703065
$ws.Range("B42").Value = 2.04186140432192
$ws.Range("C42").Value = 2.267789808285698
$ws.Range("D42").Value = 0.9905739112346479
$ws.Range("E42").Value = 1.30360487636527
$ws.Range("F42").Value = 0.8011897625074906
$ws.Range("G42").Value = 1.452823009916315
$ws.Range("H42").Value = 1.259923390874334
$ws.Range("I42").Value = 0.7491819422894392
$ws.Range("J42").Value = 0.9545826748866373
$ws.Range("K42").Value = 0.5480967329233264
$ws.Range("B43").Value = 0.2259284039637781
$ws.Range("C43").Value = -1.051287493087272
$ws.Range("D43").Value = -0.7382565279566502
$ws.Range("E43").Value = -1.240671641814429
$ws.Range("F43").Value = -0.589038394405605
$ws.Range("G43").Value = -0.781938013447586
$ws.Range("H43").Value = -1.292679462032481
$ws.Range("I43").Value = -1.087278729435283
$ws.Range("J43").Value = -1.493764671398594
$ws.Range("K43").Value = -1.058660514130252
$ws.Range("B44").Value = -1.27721589705105
$ws.Range("C44").Value = -0.9641849319204283
$ws.Range("D44").Value = -1.466600045778208
$ws.Range("E44").Value = -0.8149667983693831
$ws.Range("F44").Value = -1.007866417411364
$ws.Range("G44").Value = -1.518607865996259
$ws.Range("H44").Value = -1.313207133399061
$ws.Range("I44").Value = -1.719693075362372
$ws.Range("J44").Value = -1.28458891809403
$ws.Range("B45").Value = 0.313030965130622
$ws.Range("C45").Value = -0.1893841487271573
$ws.Range("D45").Value = 0.4622490986816672
$ws.Range("E45").Value = 0.2693494796396863
$ws.Range("F45").Value = -0.2413919689452087
$ws.Range("G45").Value = -0.0359912363480106
$ws.Range("H45").Value = -0.4424771783113215
$ws.Range("I45").Value = -0.007373021042980099
$ws.Range("B46").Value = -0.5024151138577793
$ws.Range("C46").Value = 0.1492181335510452
$ws.Range("D46").Value = -0.04368148549093578
$ws.Range("E46").Value = -0.5544229340758307
$ws.Range("F46").Value = -0.3490222014786326
$ws.Range("G46").Value = -0.7555081434419435
$ws.Range("H46").Value = -0.3204039861736021
$ws.Range("B47").Value = 0.6516332474088244
$ws.Range("C47").Value = 0.4587336283668435
$ws.Range("D47").Value = -0.05200782021805142
$ws.Range("E47").Value = 0.1533929123791467
$ws.Range("F47").Value = -0.2530930295841642
$ws.Range("G47").Value = 0.1820111276841772
$ws.Range("B48").Value = -0.192899619041981
$ws.Range("C48").Value = -0.7036410676268758
$ws.Range("D48").Value = -0.4982403350296778
$ws.Range("E48").Value = -0.9047262769929887
$ws.Range("F48").Value = -0.4696221197246473
$ws.Range("B49").Value = -0.5107414485848949
$ws.Range("C49").Value = -0.3053407159876969
$ws.Range("D49").Value = -0.7118266579510077
$ws.Range("E49").Value = -0.2767225006826664
$ws.Range("B50").Value = 0.2054007325971981
$ws.Range("C50").Value = -0.2010852093661128
$ws.Range("D50").Value = 0.2340189479022286
$ws.Range("B51").Value = -0.4064859419633109
$ws.Range("C51").Value = 0.0286182153050305
$ws.Range("B52").Value = 0.4351041572683414
